# Refresh Universalis market-price snapshot + derived profit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 557564.1
$ws.Range("I40").Value = 2250
$ws.Range("J40").Value = 1001815.4
$ws.Range("K40").Value = 2250
$ws.Range("L40").Value = 1001815.4
$ws.Range("M40").Value = -2075
$ws.Range("N40").Value = -1002165.4

# Row 127
$ws.Range("H127").Value = 2849.6086
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 2849.6086
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 8548.825800000001
$ws.Range("N127").Value = -18468.8258
$ws.Range("M127").ClearContents()

# Row 135
$ws.Range("H135").Value = 5814896.5
$ws.Range("I135").Value = 7353589
$ws.Range("J135").Value = 2059.4443
$ws.Range("K135").Value = 66182301
$ws.Range("L135").Value = 18534.9987
$ws.Range("M135").Value = -66179766
$ws.Range("N135").Value = -23604.9987

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22228414
$ws.Range("I32").Value = 4966.7075
$ws.Range("J32").Value = 250018750
$ws.Range("K32").Value = 4966.7075
$ws.Range("L32").Value = 250018750
$ws.Range("M32").Value = -4679.7075
$ws.Range("N32").Value = -250019324

# Row 60
$ws.Range("H60").Value = 7388
$ws.Range("I60").Value = 980
$ws.Range("J60").Value = 17000
$ws.Range("K60").Value = 980
$ws.Range("L60").Value = 17000
$ws.Range("M60").Value = -247
$ws.Range("N60").Value = -18466

# Row 61
$ws.Range("H61").Value = 13890235
$ws.Range("I61").Value = 13890235
$ws.Range("K61").Value = 13890235
$ws.Range("M61").Value = -13890023

# Row 74
$ws.Range("H74").Value = 1100.2034
$ws.Range("I74").Value = 1041.3489
$ws.Range("J74").Value = 1258.375
$ws.Range("K74").Value = 1041.3489
$ws.Range("L74").Value = 1258.375
$ws.Range("M74").Value = -167.3489
$ws.Range("N74").Value = -3006.375

# Row 77
$ws.Range("H77").Value = 1100.2034
$ws.Range("I77").Value = 1041.3489
$ws.Range("J77").Value = 1258.375
$ws.Range("K77").Value = 5206.7445
$ws.Range("L77").Value = 6291.875
$ws.Range("M77").Value = -838.7444999999998
$ws.Range("N77").Value = -15027.875

# Row 122
$ws.Range("H122").Value = 1426.2222
$ws.Range("I122").Value = 1480.8572
$ws.Range("K122").Value = 4442.571599999999
$ws.Range("M122").Value = -1992.571599999999

# Row 136
$ws.Range("H136").Value = 13890235
$ws.Range("I136").Value = 13890235
$ws.Range("K136").Value = 41670705
$ws.Range("M136").Value = -41668155

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 333334980
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 333334980
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 333334980
$ws.Range("N105").Value = -333338474
$ws.Range("M105").ClearContents()

# Row 107
$ws.Range("H107").Value = 41668452
$ws.Range("I107").Value = 55556610
$ws.Range("K107").Value = 55556610
$ws.Range("M107").Value = -55554690

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 198.29411
$ws.Range("I7").Value = 201.33333
$ws.Range("J7").Value = 194.875
$ws.Range("K7").Value = 201.33333
$ws.Range("L7").Value = 194.875
$ws.Range("M7").Value = -88.33332999999999
$ws.Range("N7").Value = -420.875

# Row 105
$ws.Range("H105").Value = 2805.8333
$ws.Range("I105").Value = 2851.3208
$ws.Range("J105").Value = 2461.4285
$ws.Range("K105").Value = 2851.3208
$ws.Range("L105").Value = 2461.4285
$ws.Range("M105").Value = -1104.3208
$ws.Range("N105").Value = -5955.4285

# Row 117
$ws.Range("H117").Value = 41000
$ws.Range("J117").Value = 41000
$ws.Range("L117").Value = 41000
$ws.Range("N117").Value = -50178

# Row 122
$ws.Range("H122").Value = 1615.2307
$ws.Range("I122").Value = 2649.3333
$ws.Range("J122").Value = 728.8570999999999
$ws.Range("K122").Value = 7947.999899999999
$ws.Range("L122").Value = 2186.5713
$ws.Range("M122").Value = -5497.999899999999
$ws.Range("N122").Value = -7086.5713

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 34317276
$ws.Range("I5").Value = 48611450
$ws.Range("J5").Value = 11260
$ws.Range("K5").Value = 145834350
$ws.Range("L5").Value = 33780
$ws.Range("M5").Value = -145834238
$ws.Range("N5").Value = -34004

# Row 8
$ws.Range("H8").Value = 64.666664
$ws.Range("I8").Value = 64.666664
$ws.Range("K8").Value = 193.999992
$ws.Range("M8").Value = -54.99999199999999

# Row 121
$ws.Range("H121").Value = 707.4666999999999
$ws.Range("I121").Value = 386.75
$ws.Range("J121").Value = 824.0909
$ws.Range("K121").Value = 1160.25
$ws.Range("L121").Value = 2472.2727
$ws.Range("M121").Value = 149.75
$ws.Range("N121").Value = -5092.2727

# Row 122
$ws.Range("H122").Value = 11164701
$ws.Range("I122").Value = 44643250
$ws.Range("J122").Value = 5184.952
$ws.Range("K122").Value = 401789250
$ws.Range("L122").Value = 46664.568
$ws.Range("M122").Value = -401786800
$ws.Range("N122").Value = -51564.568

# Row 126
$ws.Range("H126").Value = 83335480
$ws.Range("I126").Value = 125001560
$ws.Range("K126").Value = 375004680
$ws.Range("M126").Value = -374999740

# Row 135
$ws.Range("H135").Value = 34317276
$ws.Range("I135").Value = 48611450
$ws.Range("J135").Value = 11260
$ws.Range("K135").Value = 437503050
$ws.Range("L135").Value = 101340
$ws.Range("M135").Value = -437500515
$ws.Range("N135").Value = -106410

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 832
$ws.Range("I16").Value = 894.2857
$ws.Range("J16").Value = 723
$ws.Range("K16").Value = 894.2857
$ws.Range("L16").Value = 723
$ws.Range("M16").Value = -724.2857
$ws.Range("N16").Value = -1063

# Row 136
$ws.Range("H136").Value = 53573692
$ws.Range("I136").Value = 40818148
$ws.Range("J136").Value = 83336630
$ws.Range("K136").Value = 122454444
$ws.Range("L136").Value = 250009890
$ws.Range("M136").Value = -122451894
$ws.Range("N136").Value = -250014990

# Row 141
$ws.Range("H141").Value = 52325
$ws.Range("J141").Value = 52325
$ws.Range("L141").Value = 52325
$ws.Range("N141").Value = -62685

$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 24993.334
$ws.Range("J69").Value = 24993.334
$ws.Range("L69").Value = 24993.334
$ws.Range("N69").Value = -26491.334

# Row 72
$ws.Range("H72").Value = 24993.334
$ws.Range("J72").Value = 24993.334
$ws.Range("L72").Value = 74980.00199999999
$ws.Range("N72").Value = -82468.00199999999
